$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price column cells keep their original text formatting (avoid numeric auto-conversion)
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply updated cell values from the source diff
$ws.Range("D2").Value = '40.960.90'
$ws.Range("E2").Value = '  -4.19%  '
$ws.Range("D3").Value = '2.446.09'
$ws.Range("E3").Value = '  -3.65%  '
$ws.Range("D4").Value = '0.998'
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("D5").Value = '309.91'
$ws.Range("E5").Value = '  +0.28%  '
$ws.Range("D6").Value = '94.16'
$ws.Range("E6").Value = '  -7.47%  '
$ws.Range("E7").Value = '  -4.39%  '
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("D9").Value = '0.502'
$ws.Range("E9").Value = '  -5.06%  '
$ws.Range("D10").Value = '33.32'
$ws.Range("E10").Value = '  -8.42%  '
$ws.Range("D11").Value = '0.0779'
$ws.Range("E11").Value = '  -3.24%  '
$ws.Range("D12").Value = '0.108'
$ws.Range("E12").Value = '  -0.66%  '
$ws.Range("D13").Value = '6.94'
$ws.Range("E13").Value = '  -5.80%  '
$ws.Range("D14").Value = '2.810.52'
$ws.Range("E14").Value = '  -4.03%  '
$ws.Range("D15").Value = '2.440.08'
$ws.Range("E15").Value = '  -6.34%  '
$ws.Range("D16").Value = '14.37'
$ws.Range("E16").Value = '  -9.22%  '
$ws.Range("D17").Value = '0.782'
$ws.Range("E17").Value = '  -4.09%  '
$ws.Range("D18").Value = '40.934.15'
$ws.Range("E18").Value = '  -4.20%  '
$ws.Range("D19").Value = '6.32'
$ws.Range("E19").Value = '  -6.73%  '
$ws.Range("D20").Value = '0.0₃0911'
$ws.Range("E20").Value = '  -4.48%  '
$ws.Range("D21").Value = '11.48'
$ws.Range("E21").Value = '  -6.62%  '
$ws.Range("D22").Value = '67.11'
$ws.Range("E22").Value = '  -3.22%  '
$ws.Range("D23").Value = '236.62'
$ws.Range("E23").Value = '  -3.20%  '
$ws.Range("D24").Value = '2.76'
$ws.Range("E24").Value = '  -4.64%  '
$ws.Range("D25").Value = '1.92'
$ws.Range("E25").Value = '  -6.19%  '
$ws.Range("E26").Value = '  +0.30%  '
$ws.Range("D27").Value = '24.48'
$ws.Range("E27").Value = '  -6.35%  '
$ws.Range("D28").Value = '2.23'
$ws.Range("E28").Value = '  -4.05%  '
$ws.Range("D29").Value = '9.65'
$ws.Range("E29").Value = '  -5.47%  '
$ws.Range("D30").Value = '36.00'
$ws.Range("E30").Value = '  -8.78%  '
$ws.Range("D31").Value = '152.89'
$ws.Range("E31").Value = '  -1.95%  '
$ws.Range("D32").Value = '5.57'
$ws.Range("E32").Value = '  -3.98%  '
$ws.Range("E33").Value = '  -1.04%  '
$ws.Range("B34").Value = 'Hedera'
$ws.Range("C34").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D34").Value = '0.0749'
$ws.Range("E34").Value = '  -5.75%  '
$ws.Range("B35").Value = 'ApeXProtocol'
$ws.Range("C35").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D35").Value = '2.52'
$ws.Range("E35").Value = '  -9.02%  '
$ws.Range("D36").Value = '3.01'
$ws.Range("E36").Value = '  -5.82%  '
$ws.Range("D37").Value = '17.21'
$ws.Range("E37").Value = '  -6.52%  '
$ws.Range("D38").Value = '1.89'
$ws.Range("E38").Value = '  -7.66%  '
$ws.Range("B39").Value = 'Kaspa'
$ws.Range("C39").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D39").Value = '0.103'
$ws.Range("E39").Value = '  -8.79%  '
$ws.Range("B40").Value = 'Stellar'
$ws.Range("C40").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D40").Value = '0.113'
$ws.Range("E40").Value = '  -4.75%  '
$ws.Range("D41").Value = '4.15'
$ws.Range("E41").Value = '  -4.25%  '
$ws.Range("D42").Value = '21.15'
$ws.Range("E42").Value = '  -5.53%  '
$ws.Range("E43").Value = '  -0.03%  '
$ws.Range("D44").Value = '1.957.21'
$ws.Range("E44").Value = '  -0.99%  '
$ws.Range("D45").Value = '0.0283'
$ws.Range("E45").Value = '  -5.60%  '
$ws.Range("D46").Value = '3.03'
$ws.Range("E46").Value = '  -9.05%  '
$ws.Range("D47").Value = '8.66'
$ws.Range("E47").Value = '  -2.75%  '
$ws.Range("D48").Value = '76.63'
$ws.Range("E48").Value = '  -5.24%  '
$ws.Range("D49").Value = '96.90'
$ws.Range("E49").Value = '  -4.19%  '
$ws.Range("D50").Value = '68.85'
$ws.Range("E50").Value = '  -5.47%  '
$ws.Range("D51").Value = '0.179'
$ws.Range("E51").Value = '  -7.30%  '
